$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2379.1667
$ws.Range("I40").Value = 2315
$ws.Range("J40").Value = 2700
$ws.Range("K40").Value = 2315
$ws.Range("L40").Value = 2700
$ws.Range("M40").Value = -2140
$ws.Range("N40").Value = -3050

$ws.Range("H132").Value = 251145.84
$ws.Range("I132").Value = 304895.84
$ws.Range("J132").Value = 12256.889
$ws.Range("K132").Value = 914687.52
$ws.Range("L132").Value = 36770.667
$ws.Range("M132").Value = -912157.52
$ws.Range("N132").Value = -41830.667

$ws.Range("H137").Value = 31251144
$ws.Range("I137").Value = 43479016
$ws.Range("J137").Value = 2140.889
$ws.Range("K137").Value = 130437048
$ws.Range("L137").Value = 6422.667
$ws.Range("M137").Value = -130434498
$ws.Range("N137").Value = -11522.667

$ws.Range("H138").Value = 1696.4445
$ws.Range("I138").Value = 992.13043
$ws.Range("J138").Value = 2307.7358
$ws.Range("K138").Value = 2976.39129
$ws.Range("L138").Value = 6923.207399999999
$ws.Range("M138").Value = 2163.60871
$ws.Range("N138").Value = -17203.2074

$ws.Range("H141").Value = 2077.4568
$ws.Range("I141").Value = 1245.9155
$ws.Range("J141").Value = 7981.4
$ws.Range("K141").Value = 3737.7465
$ws.Range("L141").Value = 23944.2
$ws.Range("M141").Value = 1442.2535
$ws.Range("N141").Value = -34304.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14723.346
$ws.Range("I32").Value = 2789.081
$ws.Range("J32").Value = 103036.9
$ws.Range("K32").Value = 2789.081
$ws.Range("L32").Value = 103036.9
$ws.Range("M32").Value = -2502.081
$ws.Range("N32").Value = -103610.9

$ws.Range("H41").Value = 26500
$ws.Range("I41").Value = 3000
$ws.Range("J41").Value = 50000
$ws.Range("K41").Value = 3000
$ws.Range("L41").Value = 50000
$ws.Range("M41").Value = -2586
$ws.Range("N41").Value = -50828

$ws.Range("H74").Value = 8035.409
$ws.Range("I74").Value = 2341.3572
$ws.Range("J74").Value = 18000
$ws.Range("K74").Value = 2341.3572
$ws.Range("L74").Value = 18000
$ws.Range("M74").Value = -1467.3572
$ws.Range("N74").Value = -19748

$ws.Range("H77").Value = 8035.409
$ws.Range("I77").Value = 2341.3572
$ws.Range("J77").Value = 18000
$ws.Range("K77").Value = 11706.786
$ws.Range("L77").Value = 90000
$ws.Range("M77").Value = -7338.786
$ws.Range("N77").Value = -98736

$ws.Range("H132").Value = 4597.5293
$ws.Range("I132").Value = 4577.778
$ws.Range("J132").Value = 4619.75
$ws.Range("K132").Value = 13733.334
$ws.Range("L132").Value = 13859.25
$ws.Range("M132").Value = -11203.334
$ws.Range("N132").Value = -18919.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H49").Value = 5300
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 5300
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 5300
$ws.Range("N49").Value = -5778

$ws.Range("H134").Value = 20835694
$ws.Range("I134").Value = 27779796
$ws.Range("J134").Value = 3390.8333
$ws.Range("K134").Value = 83339388
$ws.Range("L134").Value = 10172.4999
$ws.Range("M134").Value = -83336853
$ws.Range("N134").Value = -15242.4999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 43150.918
$ws.Range("I16").Value = 63540.125
$ws.Range("J16").Value = 2372.5
$ws.Range("K16").Value = 63540.125
$ws.Range("L16").Value = 2372.5
$ws.Range("M16").Value = -63253.125
$ws.Range("N16").Value = -2946.5

$ws.Range("H31").Value = 2244.111
$ws.Range("I31").Value = 1257.625
$ws.Range("J31").Value = 4217.0835
$ws.Range("K31").Value = 1257.625
$ws.Range("L31").Value = 4217.0835
$ws.Range("M31").Value = -962.625
$ws.Range("N31").Value = -4807.0835

$ws.Range("H34").Value = 2244.111
$ws.Range("I34").Value = 1257.625
$ws.Range("J34").Value = 4217.0835
$ws.Range("K34").Value = 1257.625
$ws.Range("L34").Value = 4217.0835
$ws.Range("M34").Value = -1055.625
$ws.Range("N34").Value = -4621.0835

$ws.Range("H113").Value = 43150.918
$ws.Range("I113").Value = 63540.125
$ws.Range("J113").Value = 2372.5
$ws.Range("K113").Value = 63540.125
$ws.Range("L113").Value = 2372.5
$ws.Range("M113").Value = -61370.125
$ws.Range("N113").Value = -6712.5

$ws.Range("H132").Value = 2414.8076
$ws.Range("I132").Value = 2049.2632
$ws.Range("J132").Value = 3407
$ws.Range("K132").Value = 6147.7896
$ws.Range("L132").Value = 10221
$ws.Range("M132").Value = -3617.7896
$ws.Range("N132").Value = -15281

$ws.Range("H134").Value = 2128.6167
$ws.Range("I134").Value = 1466.76
$ws.Range("J134").Value = 5437.9
$ws.Range("K134").Value = 4400.28
$ws.Range("L134").Value = 16313.7
$ws.Range("M134").Value = -1865.28
$ws.Range("N134").Value = -21383.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 5362.5
$ws.Range("I62").Value = 450
$ws.Range("J62").Value = 7000
$ws.Range("K62").Value = 1350
$ws.Range("L62").Value = 21000
$ws.Range("M62").Value = -664
$ws.Range("N62").Value = -22372

$ws.Range("H65").Value = 5362.5
$ws.Range("I65").Value = 450
$ws.Range("J65").Value = 7000
$ws.Range("K65").Value = 4050
$ws.Range("L65").Value = 63000
$ws.Range("M65").Value = -618
$ws.Range("N65").Value = -69864

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4736.1665
$ws.Range("I132").Value = 5489.385
$ws.Range("J132").Value = 3846
$ws.Range("K132").Value = 16468.155
$ws.Range("L132").Value = 11538
$ws.Range("M132").Value = -13938.155
$ws.Range("N132").Value = -16598

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 20023.25
$ws.Range("I42").Value = 50000
$ws.Range("J42").Value = 10031
$ws.Range("K42").Value = 50000
$ws.Range("L42").Value = 10031
$ws.Range("M42").Value = -49437
$ws.Range("N42").Value = -11157

$ws.Range("H49").Value = 20023.25
$ws.Range("I49").Value = 50000
$ws.Range("J49").Value = 10031
$ws.Range("K49").Value = 50000
$ws.Range("L49").Value = 10031
$ws.Range("M49").Value = -49853
$ws.Range("N49").Value = -10325

$ws.Range("H61").Value = 4725.8887
$ws.Range("I61").Value = 11000
$ws.Range("J61").Value = 1588.8334
$ws.Range("K61").Value = 11000
$ws.Range("L61").Value = 1588.8334
$ws.Range("M61").Value = -10798
$ws.Range("N61").Value = -1992.8334

$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").ClearContents()
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = 0

$ws.Range("H113").Value = 4725.8887
$ws.Range("I113").Value = 11000
$ws.Range("J113").Value = 1588.8334
$ws.Range("K113").Value = 11000
$ws.Range("L113").Value = 1588.8334
$ws.Range("M113").Value = -8830
$ws.Range("N113").Value = -5928.8334

$ws.Range("H122").Value = 3900
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -16900

$ws.Range("H135").Value = 38000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 38000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 38000
$ws.Range("N135").Value = -48140

$ws.Range("H136").Value = 4497.6514
$ws.Range("I136").Value = 2420.9697
$ws.Range("J136").Value = 11350.7
$ws.Range("K136").Value = 7262.909100000001
$ws.Range("L136").Value = 34052.10000000001
$ws.Range("M136").Value = -4712.909100000001
$ws.Range("N136").Value = -39152.10000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6174442
$ws.Range("I132").Value = 8476104
$ws.Range("J132").Value = 1803.5
$ws.Range("K132").Value = 25428312
$ws.Range("L132").Value = 5410.5
$ws.Range("M132").Value = -25425782
$ws.Range("N132").Value = -10470.5

$ws.Range("H136").Value = 16110.288
$ws.Range("I136").Value = 22291.979
$ws.Range("J136").Value = 1892.4
$ws.Range("K136").Value = 66875.93700000001
$ws.Range("L136").Value = 5677.200000000001
$ws.Range("M136").Value = -64325.93700000001
$ws.Range("N136").Value = -10777.2
